# Auto-generated from the OOXML diff for cryptos.xlsx
# Updates coin prices / 1h volume percentages, and fixes the row-31/32
# and row-41/42 ordering swap (ImmutableX<->WrappedeETH, TheGraph<->FirstDigitalUSD).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.073.78"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.929.12"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'610.29"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").Value = "'171.51"
$ws.Range("E6").Value = "  +5.40%  "
$ws.Range("D7").Value = "3.925.01"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").Value = "'6.45"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  +5.77%  "
$ws.Range("D14").Value = "'38.61"
$ws.Range("E14").Value = "  +5.10%  "
$ws.Range("D15").Value = "4.602.60"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "3.923.02"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "70.161.01"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "'7.67"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "'18.47"
$ws.Range("E19").Value = "  +8.11%  "
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "'11.08"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "'497.33"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").Value = "'0.747"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("D24").Value = "'0.0000166"
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("D25").Value = "'86.04"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("D27").Value = "'12.40"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").Value = "'10.21"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'3.02"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("B31").Value = "WrappedeETH"
$ws.Range("C31").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D31").Value = "4.083.94"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.45"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").Value = "'7.87"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").Value = "'32.37"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "3.897.88"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").Value = "'6.17"
$ws.Range("E37").Value = "  +4.76%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "'3.29"
$ws.Range("E40").Value = "  +10.31%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.330"
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("D43").Value = "'2.12"
$ws.Range("E43").Value = "  +7.68%  "
$ws.Range("D44").Value = "'439.67"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'48.28"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'8.68"
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'0.000278"
$ws.Range("E48").Value = "  +23.65%  "
$ws.Range("D49").Value = "'0.0369"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("D50").Value = "'40.81"
$ws.Range("E50").Value = "  +6.52%  "
$ws.Range("D51").Value = "'143.26"
$ws.Range("E51").Value = "  +0.21%  "
